# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across the resume's
# achievements / work-experience bullet points.
#
# Color used for highlighted metrics: RGB(0x2C, 0x3E, 0x50) -> OOXML
# <w:color w:val="2C3E50"/>. Word's Font.Color property expects an OLE
# COLORREF (0x00BBGGRR), so we precompute the decimal value here.

$d = $word.ActiveDocument

$highlightColor = 5258796   # 0x00503E2C == RGB(0x2C,0x3E,0x50) as COLORREF
$plusMinus = [char]0x00B1   # "±"

function Highlight-Metric {
    param($Paragraph, $Phrase)

    $range = $Paragraph.Range.Duplicate
    $found = $range.Find.Execute($Phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Font.Bold = $true
        $range.Font.Color = $highlightColor
    }
    return $found
}

# --- "Partner - Siege Analytics" bullets (PROFESSIONAL EXPERIENCE) ---

# "Discovered systematic race coding errors ... from 23% to 64%"
$p1 = $d.Paragraphs(9)
Highlight-Metric $p1 "23%"
Highlight-Metric $p1 "64%"

# "Achieved 87% prediction accuracy ... from ±4.2% to ±2.1%"
$p2 = $d.Paragraphs(11)
Highlight-Metric $p2 "87%"
Highlight-Metric $p2 "71%"
$plusMinus42 = $plusMinus + "4.2%"
Highlight-Metric $p2 $plusMinus42
$plusMinus21 = $plusMinus + "2.1%"
Highlight-Metric $p2 $plusMinus21

# --- "Senior Analyst - Myers Research" bullet ---

# "Wrote RFP and analyzed bids from 1,200 vendors..."
$p3 = $d.Paragraphs(31)
Highlight-Metric $p3 "1,200"

# --- "Research Director - PCCC" bullet ---

# "...became the $400M Polling Consortium Database ... valued at $1B+"
$p4 = $d.Paragraphs(46)
Highlight-Metric $p4 "$400M"
Highlight-Metric $p4 "$1B"

# --- KEY ACHIEVEMENTS AND IMPACT bullets ---

# "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
$p5 = $d.Paragraphs(63)
Highlight-Metric $p5 "73.5%"
Highlight-Metric $p5 "$4.7M"

# "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
$p6 = $d.Paragraphs(65)
Highlight-Metric $p6 "87%"
Highlight-Metric $p6 "71%"

Write-Output "Metrics highlighted."
